# Append a new submission row to the "JSS 3B" sheet, mirroring the
# existing data layout (Timestamp, Full Name, Admission No, AI Score).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

$ws.Range("A7").Value = "2026-02-09 13:31:48"
$ws.Range("B7").Value = "USMAN BABA SHEHU"

# "10" looks numeric, but the source data keeps it as text (matching the
# existing text-typed "Number 3" / "Number 35" entries in column C).
# Force text interpretation, write it, then drop the explicit number
# format so the cell falls back to the sheet's default (unstyled) cell.
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "10"
$ws.Range("C7").ClearFormats()

$ws.Range("D7").Value = 7
